$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "27.783.84"
Set-TextValue "D3" "1.853.79"
Set-TextValue "D4" "1.014"
Set-TextValue "E4" "  -2.22%  "
Set-TextValue "D5" "320.61"
Set-TextValue "E5" "  -1.34%  "
Set-TextValue "D6" "1.013"
Set-TextValue "E6" "  -2.02%  "
Set-TextValue "D7" "0.4320"
Set-TextValue "E7" "  -2.14%  "
Set-TextValue "D8" "0.3773"
Set-TextValue "E8" "  -0.59%  "
Set-TextValue "D9" "0.07408"
Set-TextValue "E9" "  -0.73%  "
Set-TextValue "D10" "0.8854"
Set-TextValue "E10" "  +0.09%  "
Set-TextValue "D11" "21.76"
Set-TextValue "E11" "  -0.13%  "
Set-TextValue "D12" "1.870.75"
Set-TextValue "E12" "  +0.03%  "
Set-TextValue "D13" "6.771"
Set-TextValue "E13" "  +0.26%  "
Set-TextValue "D14" "5.491"
Set-TextValue "E14" "  -1.19%  "
Set-TextValue "D15" "0.07117"
Set-TextValue "E15" "  -1.54%  "
Set-TextValue "D16" "88.47"
Set-TextValue "E16" "  +5.67%  "
Set-TextValue "E17" "  -1.97%  "
Set-TextValue "D18" "0.000009050"
Set-TextValue "E18" "  -0.97%  "
Set-TextValue "E19" "  -2.04%  "
Set-TextValue "E20" "  -0.06%  "
Set-TextValue "D21" "27.792.50"
Set-TextValue "E21" "  -0.02%  "
Set-TextValue "D22" "5.286"
Set-TextValue "E22" "  -0.59%  "
Set-TextValue "D23" "11.23"
Set-TextValue "E23" "  -1.53%  "
Set-TextValue "D24" "2.082.34"
Set-TextValue "E24" "  -0.54%  "
Set-TextValue "D25" "2.038"
Set-TextValue "E25" "  +3.56%  "
Set-TextValue "D26" "156.36"
Set-TextValue "E26" "  -1.26%  "
Set-TextValue "E27" "  -0.97%  "
Set-TextValue "D28" "2.118"
Set-TextValue "E28" "  +6.39%  "
Set-TextValue "D29" "5.457"
Set-TextValue "E29" "  +2.55%  "
Set-TextValue "D30" "121.58"
Set-TextValue "E30" "  +3.22%  "
Set-TextValue "D31" "0.08974"
Set-TextValue "E31" "  -1.50%  "
Set-TextValue "D32" "1.242"
Set-TextValue "E32" "  +1.92%  "
Set-TextValue "D33" "0.7842"
Set-TextValue "E33" "  +1.08%  "
Set-TextValue "D34" "4.594"
Set-TextValue "E34" "  +0.13%  "
Set-TextValue "D35" "2.925"
Set-TextValue "E35" "  -3.91%  "
Set-TextValue "D36" "1.149"
Set-TextValue "E36" "  -1.55%  "
Set-TextValue "D37" "1.014"
Set-TextValue "E37" "  -2.12%  "
Set-TextValue "D38" "0.05342"
Set-TextValue "E38" "  -0.16%  "
Set-TextValue "E39" "  -0.87%  "
Set-TextValue "D40" "2.873"
Set-TextValue "E40" "  +1.02%  "
Set-TextValue "D41" "7.120"
Set-TextValue "E41" "  +3.25%  "
Set-TextValue "D42" "0.5214"
Set-TextValue "E42" "  +0.18%  "
Set-TextValue "D43" "0.1689"
Set-TextValue "E43" "  -0.36%  "
Set-TextValue "D44" "8.994"
Set-TextValue "E44" "  +3.06%  "
Set-TextValue "D45" "111.14"
Set-TextValue "E45" "  +1.32%  "
Set-TextValue "D46" "10.77"
Set-TextValue "E46" "  +0.48%  "
Set-TextValue "D47" "1.722"
Set-TextValue "E47" "  -0.03%  "
Set-TextValue "D48" "0.4761"
Set-TextValue "E48" "  +1.10%  "
Set-TextValue "D49" "0.06516"
Set-TextValue "E49" "  +1.17%  "
Set-TextValue "D50" "1.014"
Set-TextValue "E50" "  -2.25%  "
Set-TextValue "D51" "1.909"
Set-TextValue "E51" "  +1.31%  "
